$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:G values per row (row index => values for columns B,C,D,E,F,G)
$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    3  = @(3.230985683306322, 3099.503889238888, 3.900430680208489, 8.660232485948974, 1, 3115.295538088352)
    4  = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1, 2.290389397800092)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    6  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 3.645393585217082)
    7  = @(0.127881588408715,  1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 2.449980674824537)
    8  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    9  = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 645.3272768299601, 1, 676.4434635367506)
    10 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 8.660232485948974, 0, 38.00504557982321)
    11 = @(0.3048080303191223, 114.8270160096505, 0.8054896365839992, 8.660232485948974, 1, 124.5975461625026)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
